$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.837.33'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').Value = '3.281.78'
$ws.Range('E3').Value = '  +0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '586.26'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '178.32'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.62%  '
$ws.Range('E7').Value = '  +1.84%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '3.284.03'
$ws.Range('E9').Value = '  +0.52%  '
$ws.Range('E10').Value = '  -0.51%  '
$ws.Range('E11').Value = '  +2.17%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.400'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('D13').Value = '3.857.79'
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('E14').Value = '  -2.71%  '
$ws.Range('D15').Value = '65.952.75'
$ws.Range('E15').Value = '  +0.08%  '
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('E17').Value = '  +0.23%  '
$ws.Range('D18').Value = '3.278.69'
$ws.Range('E18').Value = '  +1.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '420.87'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.47'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.02'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.25'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.74%  '
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.09'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.52%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.66'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.38%  '
$ws.Range('E26').Value = '  +5.98%  '
$ws.Range('E27').Value = '  +0.26%  '
$ws.Range('E28').Value = '  +0.92%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.38'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +6.02%  '
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('E31').Value = '  -0.55%  '
$ws.Range('E32').Value = '  -0.28%  '
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.13'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.56'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.26%  '
$ws.Range('E36').Value = '  +0.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '158.07'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.19%  '
$ws.Range('E38').Value = '  -0.21%  '
$ws.Range('D39').Value = '2.850.61'
$ws.Range('E39').Value = '  +3.34%  '
$ws.Range('E40').Value = '  -0.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '26.16'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.90%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.33'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.54%  '
$ws.Range('E43').Value = '  -3.59%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '39.58'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.66%  '
$ws.Range('E45').Value = '  -2.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0636'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.26'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.57%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '310.86'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.11%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '22.82'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0268'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.55%  '
$ws.Range('E51').Value = '  -0.24%  '
